# Atualização de bases das ligas, do dia: 11-06-2024 às 21:19
# Swap the contents of rows 137/138 and rows 236/237 (columns B..AD),
# keeping column A (the sequential id) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($rowA, $rowB) {
    $rangeA = $ws.Range("B$rowA`:AD$rowA")
    $rangeB = $ws.Range("B$rowB`:AD$rowB")

    $valuesA = $rangeA.Value()
    $valuesB = $rangeB.Value()

    $rangeA.Value = $valuesB
    $rangeB.Value = $valuesA
}

Swap-Rows 137 138
Swap-Rows 236 237
